# Applies the cryptos-list price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.372.65'
$ws.Range('E2').Value = '  +3.55%  '
$ws.Range('D3').Value = '2.317.79'
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('E4').Value = '  -0.01%  '
$style_D5 = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '545.12'
$ws.Range('D5').Style = $style_D5
$ws.Range('E5').Value = '  +1.56%  '
$style_D6 = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.91'
$ws.Range('D6').Style = $style_D6
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -1.80%  '
$ws.Range('D9').Value = '2.313.76'
$ws.Range('E9').Value = '  +1.19%  '
$ws.Range('E10').Value = '  +0.23%  '
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('E13').Value = '  +0.26%  '
$ws.Range('E14').Value = '  -0.29%  '
$ws.Range('D15').Value = '60.297.54'
$ws.Range('E15').Value = '  +3.51%  '
$ws.Range('D16').Value = '2.730.98'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('D18').Value = '2.317.80'
$ws.Range('E18').Value = '  +2.29%  '
$style_D19 = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.56'
$ws.Range('D19').Style = $style_D19
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('E20').Value = '  -1.55%  '
$style_D21 = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '313.86'
$ws.Range('D21').Style = $style_D21
$ws.Range('E21').Value = '  -0.22%  '
$style_D22 = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.64'
$ws.Range('D22').Style = $style_D22
$ws.Range('E22').Value = '  +1.30%  '
$style_D23 = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('D23').Style = $style_D23
$ws.Range('E23').Value = '  -0.31%  '
$style_D24 = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.81'
$ws.Range('D24').Style = $style_D24
$ws.Range('E24').Value = '  +0.80%  '
$ws.Range('E25').Value = '  +1.37%  '
$ws.Range('E26').Value = '  +0.04%  '
$style_D27 = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.82'
$ws.Range('D27').Style = $style_D27
$ws.Range('E27').Value = '  -2.16%  '
$style_D28 = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.35'
$ws.Range('D28').Style = $style_D28
$ws.Range('E28').Value = '  +3.68%  '
$style_D29 = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '173.28'
$ws.Range('D29').Style = $style_D29
$ws.Range('E29').Value = '  +1.27%  '
$style_D30 = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.18'
$ws.Range('D30').Style = $style_D30
$ws.Range('E30').Value = '  +8.59%  '
$style_D31 = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.72'
$ws.Range('D31').Style = $style_D31
$ws.Range('E31').Value = '  +0.84%  '
$ws.Range('D32').Value = '0.0₃0728'
$ws.Range('E32').Value = '  +0.32%  '
$style_D33 = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.91'
$ws.Range('D33').Style = $style_D33
$ws.Range('E33').Value = '  +1.53%  '
$style_D34 = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.37'
$ws.Range('D34').Style = $style_D34
$ws.Range('E34').Value = '  +10.34%  '
$ws.Range('E35').Value = '  -0.43%  '
$style_D37 = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.82'
$ws.Range('D37').Style = $style_D37
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('E39').Value = '  +2.75%  '
$style_D40 = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '324.45'
$ws.Range('D40').Style = $style_D40
$ws.Range('E40').Value = '  +11.81%  '
$style_D41 = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '37.90'
$ws.Range('D41').Style = $style_D41
$ws.Range('E41').Value = '  -1.21%  '
$style_D42 = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.52'
$ws.Range('D42').Style = $style_D42
$ws.Range('E42').Value = '  +1.35%  '
$style_D43 = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '137.70'
$ws.Range('D43').Style = $style_D43
$ws.Range('E43').Value = '  -1.83%  '
$style_D44 = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.48'
$ws.Range('D44').Style = $style_D44
$ws.Range('E44').Value = '  +0.82%  '
$ws.Range('E45').Value = '  -1.27%  '
$style_D46 = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.16'
$ws.Range('D46').Style = $style_D46
$ws.Range('E46').Value = '  +4.91%  '
$style_D47 = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0494'
$ws.Range('D47').Style = $style_D47
$ws.Range('E47').Value = '  -0.24%  '
$ws.Range('E48').Value = '  +0.67%  '
$style_D49 = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0212'
$ws.Range('D49').Style = $style_D49
$ws.Range('E49').Value = '  +0.71%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0211'
$ws.Range('E50').Value = '  +16.68%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$style_D51 = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.03'
$ws.Range('D51').Style = $style_D51
$ws.Range('E51').Value = '  +0.78%  '
